$wb = $excel.ActiveWorkbook

# --- Sheet "Examples part 1": rename the common keyword "xlsx*" to "tbs:*" ---
$ws1 = $wb.Worksheets.Item("Examples part 1")

$ws1.Range("C26").Value = "[cell2.score;block=tbs:cell;ope=tbs:num]"

$ws1.Range("C34").Value = "tbs:num"
$ws1.Range("C35").Value = "tbs:bool"
$ws1.Range("C36").Value = "tbs:date"

$ws1.Range("D36").Value = "[onshow.x_dt;ope=tbs:date]"
$ws1.Range("D35").Value = "[onshow.x_bt;ope=tbs:bool]"
$ws1.Range("D34").Value = "[onshow.x_num;ope=tbs:num]"

$ws1.Range("E20").Value = "[a.score;ope=tbs:num]"
$ws1.Range("F20").Value = "[a.score;ope=tbs:num]"

# --- Sheet "Delete me": add a new cell referenced by a workbook-level name ---
$wsDelete = $wb.Worksheets.Item("Delete me")
$wsDelete.Range("B6").Value = "And this named cell too."

$wb.Names.Add("the_named_cell", "='Delete me'!`$B`$6")

# --- Back to "Examples part 1": rename the second "Score" header to "Score again" ---
$ws1.Range("F19").Value = "Score again"
